# Update countries & provincias Spain
# - Re-rank "Arabia Saudita" ahead of "Mexico" / "Indonesia" (rows 39-41)
# - Refresh case counters for rows 18, 39, 40, 41, 77, 108, 132
# - Bump the "Datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (row 1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 14:52"

# Row 18 - Portugal
$ws.Range("B18").Value = 15987
$ws.Range("C18").Value = 515
$ws.Range("D18").Value = 266
$ws.Range("E18").Value = 15251
$ws.Range("F18").Value = 233
$ws.Range("G18").Value = 35
$ws.Range("H18").Value = 470

# Row 39 - now Arabia Saudita (was Mexico)
$ws.Range("A39").Value = "Arabia Saudita"
$ws.Range("B39").Value = 4033
$ws.Range("C39").Value = 382
$ws.Range("D39").Value = 720
$ws.Range("E39").Value = 3261
$ws.Range("F39").Value = 57
$ws.Range("G39").Value = 5
$ws.Range("H39").Value = 52

# Row 40 - now Mexico (was Indonesia)
$ws.Range("A40").Value = "Mexico"
$ws.Range("B40").Value = 3844
$ws.Range("C40").Value = 403
$ws.Range("D40").Value = 633
$ws.Range("E40").Value = 2978
$ws.Range("F40").Value = 89
$ws.Range("G40").Value = 39
$ws.Range("H40").Value = 233

# Row 41 - now Indonesia (was Arabia Saudita)
$ws.Range("A41").Value = "Indonesia"
$ws.Range("B41").Value = 3842
$ws.Range("C41").Value = 330
$ws.Range("D41").Value = 286
$ws.Range("E41").Value = 3229
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 21
$ws.Range("H41").Value = 327

# Row 77 - Republica de Macedonia
$ws.Range("E77").Value = 685
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = 34

# Row 108 - Bolivia
$ws.Range("B108").Value = 275
$ws.Range("C108").Value = 7
$ws.Range("E108").Value = 253
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 20

# Row 132 - Madagascar
$ws.Range("B132").Value = 102
$ws.Range("C132").Value = 9
$ws.Range("E132").Value = 91
